$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values: row -> D, J, K, L, M, P (weekly data refresh / reshuffle)
$data = @(
    ,@(2, 44330, 250, 2800, 3000, 2900, 2900)
    ,@(3, 44659, 250, 950, 1000, 975, 975)
    ,@(4, 44442, 240, 2300, 2500, 2400, 2400)
    ,@(5, 44540, 200, 900, 1000, 950, 950)
    ,@(6, 44349, 250, 2800, 3000, 2900, 2900)
    ,@(7, 44571, 250, 900, 1000, 950, 950)
    ,@(8, 44326, 200, 2700, 2800, 2750, 2750)
    ,@(9, 44669, 300, 950, 1000, 975, 975)
    ,@(10, 44474, 250, 2000, 2500, 2250, 2250)
    ,@(11, 44727, 270, 1500, 2000, 1750, 1750)
    ,@(12, 44523, 250, 1400, 1500, 1450, 1450)
    ,@(13, 44536, 250, 900, 1000, 950, 950)
    ,@(14, 44224, 200, 750, 800, 775, 775)
    ,@(15, 44539, 300, 900, 1000, 950, 950)
    ,@(16, 44603, 250, 2500, 3000, 2750, 2750)
    ,@(17, 44699, 300, 2000, 2500, 2250, 2250)
    ,@(18, 44532, 300, 1000, 1200, 1100, 1100)
    ,@(19, 44250, 250, 1000, 1200, 1100, 1100)
    ,@(20, 44635, 300, 1900, 2000, 1950, 1950)
    ,@(21, 44249, 200, 900, 1000, 950, 950)
    ,@(22, 44302, 200, 900, 1000, 950, 950)
    ,@(23, 44664, 250, 1300, 1500, 1400, 1400)
    ,@(24, 44435, 300, 2300, 2500, 2400, 2400)
    ,@(25, 44498, 270, 2000, 2300, 2150, 2150)
    ,@(26, 44376, 270, 2400, 2500, 2437, 2437)
    ,@(27, 44274, 250, 1000, 1200, 1100, 1100)
    ,@(28, 44432, 300, 2300, 2500, 2400, 2400)
    ,@(29, 44645, 300, 1200, 1500, 1350, 1350)
    ,@(30, 44260, 250, 900, 1000, 950, 950)
    ,@(31, 44365, 250, 2400, 2500, 2450, 2450)
    ,@(32, 44313, 250, 900, 1000, 950, 950)
    ,@(33, 44616, 200, 2500, 3000, 2750, 2750)
    ,@(34, 44292, 250, 1800, 2000, 1900, 1900)
    ,@(35, 44417, 250, 4000, 4500, 4250, 4250)
    ,@(36, 44280, 250, 1400, 1500, 1450, 1450)
    ,@(37, 44628, 250, 2500, 3000, 2750, 2750)
    ,@(38, 44494, 200, 2400, 2500, 2450, 2450)
    ,@(39, 44362, 250, 2800, 3000, 2900, 2900)
    ,@(40, 44747, 250, 2000, 2500, 2250, 2250)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}
